$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text (A1, B1, C1): link_testSuiteLinks_* -> link_testSuiteActions_*
$ws.Range("A1").Value = "link_testSuiteActions_executions_id"
$ws.Range("B1").Value = "link_testSuiteActions_project_id"
$ws.Range("C1").Value = "link_testSuiteActions_team_id"

# Update column widths: A 35->37, B 32->34, C 29->31
# (ColumnWidth values offset by -0.85 to compensate for Excel's pixel-rounding
#  so the resulting stored XML "width" attribute lands exactly on the integer target)
$ws.Columns.Item(1).ColumnWidth = 36.15
$ws.Columns.Item(2).ColumnWidth = 33.15
$ws.Columns.Item(3).ColumnWidth = 30.15
